# Updated cryptos list row values (Price / Volume(1h)) to match the refreshed
# market snapshot, including the Hedera/InternetComputer/VeChain/TrustWalletToken
# row re-ordering that came with this refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '20.015.55'
$ws.Range("E2").Value = '  -8.02%  '
$ws.Range("D3").Value = '1.424.46'
$ws.Range("E3").Value = '  -7.52%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.003'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '273.54'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.76%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3755'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3088'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '40.16'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -7.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.010'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06579'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -8.52%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.004'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.27%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.371'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.12'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -7.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.157'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.88%  '
$ws.Range("D16").Value = '1.432.03'
$ws.Range("E16").Value = '  -7.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001011'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -8.26%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.05812'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -11.79%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '75.39'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -9.60%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.666'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -7.85%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.44'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.332'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.30%  '
$ws.Range("D25").Value = '20.038.05'
$ws.Range("E25").Value = '  -7.95%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.279'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.41%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '138.77'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.83'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -8.49%  '
$ws.Range("D29").Value = '1.590.65'
$ws.Range("E29").Value = '  -7.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '108.98'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.44%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.897'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -19.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8957'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.399'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -8.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.07782'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.383'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.38%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '11.31'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.002'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.12%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05699'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.95%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.742'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.48%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.1912'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.49%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02019'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -8.56%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.110'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.33%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.279'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -12.93%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5316'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.539'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.25'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5126'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.778'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '109.62'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -7.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.049'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.002'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.11%  '
